$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'44.171.76"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.83%  '

$ws.Range('D3').Value = "'2.251.04"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.64%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = "'318.67"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.27%  '

$ws.Range('D6').Value = "'101.36"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.04%  '

$ws.Range('E7').Value = '  -1.44%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').Value = "'0.546"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('D10').Value = "'36.97"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.54%  '

$ws.Range('D11').Value = "'0.0828"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.19%  '

$ws.Range('D12').Value = "'7.54"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.19%  '

$ws.Range('E13').Value = '  -1.91%  '

$ws.Range('D14').Value = "'2.593.40"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.68%  '

$ws.Range('D15').Value = "'2.283.13"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.86%  '

$ws.Range('D16').Value = "'0.850"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.91%  '

$ws.Range('D17').Value = "'14.23"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.06%  '

$ws.Range('D18').Value = "'44.093.58"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.75%  '

$ws.Range('D19').Value = "'13.38"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.67%  '

$ws.Range('D20').Value = "'0.0₃0978"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.67%  '

$ws.Range('D21').Value = "'6.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.89%  '

$ws.Range('D22').Value = "'65.60"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.48%  '

$ws.Range('D23').Value = "'3.09"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.22%  '

$ws.Range('D24').Value = "'235.40"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.28%  '

$ws.Range('D25').Value = "'2.07"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.53%  '

$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.05%  '

$ws.Range('D27').Value = "'10.51"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.13%  '

$ws.Range('D28').Value = "'2.21"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.25%  '

$ws.Range('D29').Value = "'37.81"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.87%  '

$ws.Range('D30').Value = "'6.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.36%  '

$ws.Range('D31').Value = "'158.87"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.85%  '

$ws.Range('D32').Value = "'20.11"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.95%  '

$ws.Range('D33').Value = "'0.0847"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.00%  '

$ws.Range('E34').Value = '  -1.00%  '

$ws.Range('D35').Value = "'3.19"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.12%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = "'0.112"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.39%  '

$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'1.95"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.17%  '

$ws.Range('E38').Value = '  -1.98%  '

$ws.Range('D39').Value = "'16.14"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.80%  '

$ws.Range('E40').Value = '  -0.68%  '

$ws.Range('D41').Value = "'4.16"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.71%  '

$ws.Range('E42').Value = '  -2.29%  '

$ws.Range('E43').Value = '  +0.14%  '

$ws.Range('D44').Value = "'1.744.59"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.98%  '

$ws.Range('D45').Value = "'0.198"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.89%  '

$ws.Range('D46').Value = "'82.56"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.91%  '

$ws.Range('D47').Value = "'74.73"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.74%  '

$ws.Range('D48').Value = "'5.16"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.88%  '

$ws.Range('E49').Value = '  +4.18%  '

$ws.Range('D50').Value = "'102.78"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.73%  '

$ws.Range('D51').Value = "'57.70"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.19%  '
